# Economic Dashboard data refresh - 2026-01-07
# Shifts "Present/Lag1..Lag4" columns, updates "Latest Period" dates, and
# toggles the yellow "just refreshed" highlight on the N (date) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Highlight maintenance on column N (Latest Period date cells).
#    Rows that received fresh data this refresh get the yellow highlight;
#    rows that were highlighted from the prior refresh but were not
#    updated this time revert to the plain (no-fill) look.
#    We do this via Copy/PasteSpecial(Formats) from donor cells that
#    already carry the desired look, so existing style entries are
#    reused instead of minting new ones.
# ---------------------------------------------------------------------

# Donor with the "highlighted" (yellow) date style.
$ws.Range("N40").Copy()
$ws.Range("N5").PasteSpecial(-4122)

# Donor with the plain (no highlight) date style.
$ws.Range("N3").Copy()
$ws.Range("N13").PasteSpecial(-4122)

$ws.Range("N4").Copy()
$ws.Range("N14").PasteSpecial(-4122)

$ws.Range("N3").Copy()
$ws.Range("N51").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Data updates
# ---------------------------------------------------------------------

# Row 5 - ADP Employment (ADPMNUSNERSA)
$ws.Range("N5").Value = 45992
$ws.Range("Q5").Value = 41000
$ws.Range("R5").Value = -29000
$ws.Range("S5").Value = 47000
$ws.Range("T5").Value = -29000
$ws.Range("U5").Value = -3000

# Row 29 - 5-Year Forward Inflation Expectation (T5YIFR)
$ws.Range("N29").Value = 46028
$ws.Range("Q29").Value = 2.24
$ws.Range("R29").Value = 2.23
$ws.Range("S29").Value = 2.22
$ws.Range("T29").Value = 2.24
$ws.Range("U29").Value = 2.23

# Row 30 - 10-Year Breakeven Inflation Rate (T10YIE)
$ws.Range("N30").Value = 46028
$ws.Range("Q30").Value = 2.27
$ws.Range("R30").Value = 2.26
$ws.Range("S30").Value = 2.25
$ws.Range("T30").Value = 2.25
$ws.Range("U30").Value = 2.24

# Row 47 - Effective Federal Funds Rate (DFF)
$ws.Range("N47").Value = 46027

# Row 48 - 2-Year Treasury Yield (DGS2)
$ws.Range("N48").Value = 46027
$ws.Range("Q48").Value = 3.46
$ws.Range("R48").Value = 3.47
$ws.Range("S48").Value = 3.47
$ws.Range("T48").Value = 3.45
$ws.Range("U48").Value = 3.45

# Row 49 - 5-Year Treasury Yield (DGS5)
$ws.Range("N49").Value = 46027
$ws.Range("Q49").Value = 3.71
$ws.Range("R49").Value = 3.74
$ws.Range("S49").Value = 3.73
$ws.Range("T49").Value = 3.68
$ws.Range("U49").Value = 3.67

# Row 50 - 10-Year Treasury Yield (DGS10)
$ws.Range("N50").Value = 46027
$ws.Range("Q50").Value = 4.17
$ws.Range("R50").Value = 4.19
$ws.Range("S50").Value = 4.18
$ws.Range("T50").Value = 4.14
$ws.Range("U50").Value = 4.12

# Row 52 - Corporate bond yield
$ws.Range("N52").Value = 46027
$ws.Range("Q52").Value = 5.92
$ws.Range("R52").Value = 5.93
$ws.Range("S52").Value = 5.9
$ws.Range("T52").Value = 5.89
$ws.Range("U52").Value = 5.88
